# Update profit files after running on 2025-09-27
# Append the next day's row (date + profit) to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 41

# Write the date as literal text (leading apostrophe forces text, so the
# date-looking string "09/27/2025" isn't auto-converted into a real Excel
# date serial number), matching how the rest of the column is stored.
$ws.Range("A" + $newRow).Value = "'09/27/2025"
# Re-normalize the style so the quote-prefix/text formatting Excel applies
# when forcing text doesn't leave a stray number format on the cell - the
# rest of the date column has no explicit style either.
$ws.Range("A" + $newRow).Style = "Normal"

# Profit value for that date - a plain number.
$ws.Range("B" + $newRow).Value = 14542.59
